$d = $word.ActiveDocument
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*tipo mais comum e o mais recomendado*") {
        $target = $p
        break
    }
}

$xmlFrag = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5863641A" w14:textId="5B4D17A5" w:rsidR="00D369C6" w:rsidRDefault="00D369C6" w:rsidP="00E36DA2"><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t>Também existem tipos de servidores como, servidor de torre que é o tipo mais comum e o mais recomendado para pequenas empresas, com uma infraestrutura</w:t></w:r><w:r w:rsidR="007C4036"><w:t xml:space="preserve"> de TI menor,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>monitoramento e manutenção de forma mais simples. Um outro tipo de servidor chamado Servidor em Rack que é indicado para empresas já em crescimento, permitindo expansão, conexão de unidades de armazenamento externo do tipo NAS ou SAN e arquitetura altamente dimensionável</w:t></w:r><w:r><w:t xml:space="preserve">. E por último o servidor Infraestrutura Modular que é típico indicado para grandes empresas, que necessita de uma alta capacidade de processamento e armazenamento, onde sua principal vantagem é que os cabos de alimentação, cabos de rede e sistemas de arrefecimento são compartilhados entre os Servidores </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Blade</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Ou seja</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> existem outros tipos de servidores para tipos de capacidades diferentes.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.Range.InsertXML($xmlFrag)
